$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear row 5 (Eric's entry) completely - contents and formatting, no shifting of row 30
$ws.Range("A5:C5").Clear()

# Update selection to A5 (post-edit sheet view)
$ws.Range("A5").Select()
